# Added new param for calibration
# - Move the "x" marker in the "Check to include in analysis" column (E)
#   on the "Optimal funding scenario" sheet from Scenario 1 (row 2) to
#   Scenario 3 (row 4).
# - Update the remembered cell selections on both sheets to match where the
#   author was last working.

$wb = $excel.ActiveWorkbook

$currentExpenditure = $wb.Worksheets.Item("Current expenditure")
$optimalFunding = $wb.Worksheets.Item("Optimal funding scenario")

# --- Core data edit: re-flag which scenario should be calibrated ---
$optimalFunding.Range("E2").ClearContents()
$optimalFunding.Range("E4").Value = "x"

# --- View state: restore author's last selections / scroll position ---
$currentExpenditure.Activate()
$currentExpenditure.Range("G16").Select()

$optimalFunding.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$optimalFunding.Range("E5").Select()
